$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (before current row 2), pushing
# the existing 20 data rows down to rows 4-23.
$ws.Rows.Item(2).Resize(2).Insert()
# Remove formatting Excel copied from the header row during insert so the
# new rows stay plain (matching the rest of the data rows).
$ws.Rows.Item(2).Resize(2).ClearFormats()

# Values for the two newly inserted rows.
$ws.Cells.Item(2,1).Value = -0.1967945098876953
$ws.Cells.Item(2,2).Value = -0.0881298780441284
$ws.Cells.Item(2,3).Value = -0.4354097247123718

$ws.Cells.Item(3,1).Value = 0.3726930618286133
$ws.Cells.Item(3,2).Value = 0.3928739428520202
$ws.Cells.Item(3,3).Value = -0.1955753564834594

# Append 8 new rows after the existing data (now ending at row 23).
$ws.Cells.Item(24,1).Value = 0.6416101455688477
$ws.Cells.Item(24,2).Value = 0.1077315807342529
$ws.Cells.Item(24,3).Value = -3.572567462921143

$ws.Cells.Item(25,1).Value = 1.180892944335938
$ws.Cells.Item(25,2).Value = -0.3624088764190674
$ws.Cells.Item(25,3).Value = 1.944910764694214

$ws.Cells.Item(26,1).Value = -0.6099348068237305
$ws.Cells.Item(26,2).Value = -0.0995303392410278
$ws.Cells.Item(26,3).Value = 1.559979677200317

$ws.Cells.Item(27,1).Value = 0.7382268905639648
$ws.Cells.Item(27,2).Value = 0.5965696573257446
$ws.Cells.Item(27,3).Value = 0.3601601719856262

$ws.Cells.Item(28,1).Value = 0.1256790161132812
$ws.Cells.Item(28,2).Value = 0.4359270334243774
$ws.Cells.Item(28,3).Value = -0.5883067846298218

$ws.Cells.Item(29,1).Value = 0.1187114715576171
$ws.Cells.Item(29,2).Value = 0.2241333723068237
$ws.Cells.Item(29,3).Value = -0.9467962980270386

$ws.Cells.Item(30,1).Value = 0.1263256072998047
$ws.Cells.Item(30,2).Value = 0.5689128637313843
$ws.Cells.Item(30,3).Value = -0.7026804089546204

$ws.Cells.Item(31,1).Value = -0.1413173675537109
$ws.Cells.Item(31,2).Value = 0.4839025735855102
$ws.Cells.Item(31,3).Value = -0.0290583968162536
